$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append two new rows of feed log data to the bottom of the sheet (rows 60 and 61)
$ws.Range("A60").Value = 59
$ws.Range("B60").Value = 1
$ws.Range("C60").Value = "2024-06-16 04:14:40"
$ws.Range("D60").Value = 200
$ws.Range("E60").Value = 3

$ws.Range("A61").Value = 60
$ws.Range("B61").Value = 2
$ws.Range("C61").Value = "2024-06-16 04:14:41"
$ws.Range("D61").Value = 200
$ws.Range("E61").Value = 0
